$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("寄件件數 / Số lượng" - quantity shipped) updates per row
$qty = @{
    2  = 6
    3  = 3
    4  = 6
    6  = 5
    7  = 5
    8  = 6
    9  = 8
    10 = 2
    11 = 8
    12 = 8
    13 = 7
    14 = 4
    16 = 3
    17 = 9
    19 = 8
    20 = 10
    21 = 2
    22 = 10
    23 = 6
    25 = 6
    26 = 3
    27 = 4
    28 = 10
    29 = 9
    30 = 10
    31 = 4
    32 = 5
    33 = 7
    35 = 9
    36 = 4
    37 = 4
    38 = 4
    39 = 3
    46 = 2
}

foreach ($row in $qty.Keys) {
    $ws.Cells.Item($row, 2).Value = $qty[$row]
}

# Row 39 previously had no recipient-name (column H) value because the
# product info was turned off; set it to "None" to match column G/I.
$ws.Cells.Item(39, 8).Value = "None"
